$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Update the "education" sheet: the PhD entry's description no
#    longer says "Manuscript on review" and instead gets a proper
#    italicised species name.
# ---------------------------------------------------------------------
$edu = $wb.Worksheets.Item("education")
$full = "Title: Phylogeography and species distribucion models to define Conservation Units of the southern three-banded armadillo (Tolypeutes matacus) in Argentina."
$edu.Range("E2").Value = $full

$speciesStart = $full.IndexOf("Tolypeutes matacus") + 1
$speciesLen = "Tolypeutes matacus".Length
$afterStart = $speciesStart + $speciesLen

$italicRun = $edu.Range("E2").Characters($speciesStart, $speciesLen)
$italicRun.Font.Italic = $true
$italicRun.Font.ColorIndex = -4105

$restLen = $full.Length - ($afterStart - 1)
$restRun = $edu.Range("E2").Characters($afterStart, $restLen)
$restRun.Font.Size = 11
$restRun.Font.Name = "Calibri"
$restRun.Font.ColorIndex = -4105

$edu.Range("E2").Select()

# ---------------------------------------------------------------------
# 2) Insert a new worksheet "msSubmitted" right after "education" and
#    before "awards" (i.e. as the 2nd sheet in the workbook).
# ---------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "msSubmitted"
$newSheet.Move($wb.Worksheets.Item(2))

# Re-fetch a fresh reference by name (the old object handle goes stale
# once the sheet has been repositioned).
$ms = $wb.Worksheets.Item("msSubmitted")

# Column widths matching the source CV layout.
$ms.Columns.Item(1).ColumnWidth = 138.625
$ms.Columns.Item(2).ColumnWidth = 28.375
$ms.Columns.Item(3).ColumnWidth = 28

# Header row.
$ms.Range("A1").Value = "title"
$ms.Range("B1").Value = "journal"
$ms.Range("C1").Value = "status"
$ms.Range("A1:C1").Font.Bold = $true
$ms.Range("A1:C1").Font.Size = 12

# Data row: the manuscript currently under submission.
$ms.Range("A2").Value = "Reconstructing the distribution of chacoan biota from current and past evidence: the case of the southern three-banded armadillo Tolypeutes matacus (Desmarest, 1804)"
$ms.Range("B2").Value = "Journal of Mammalian Evolution"
$ms.Range("C2").Value = "Under second round of reviews"

# ---------------------------------------------------------------------
# 3) Make "msSubmitted" the active/selected sheet, as in the source file.
# ---------------------------------------------------------------------
$ms.Select()
